$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header tweaks
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 - "deleted" C2, updated B2/D2/E2
$ws.Range("B2").Value = 5.2025224677145037
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 11.54005608882097
$ws.Range("E2").Value = 10.303267829789519

# Row 3 - "deleted" B3, added D3, updated C3/E3
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 6.9618007501512436
$ws.Range("D3").Value = 5.8348200070450318
$ws.Range("E3").Value = 11.575310008874508

# Selection change
$ws.Range("B1:E3").Select()
